$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "If using Project Scarlett, set the active solution platform to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If using an Xbox One X|S devkit, set the active solution platform to ",
    2)

Write-Output "done"
